$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 316.22223
$ws.Range("I2").Value = 237.8
$ws.Range("J2").Value = 414.25
$ws.Range("K2").Value = 237.8
$ws.Range("L2").Value = 414.25
$ws.Range("M2").Value = -124.8
$ws.Range("N2").Value = -640.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 259.30768
$ws.Range("I12").Value = 214.63637
$ws.Range("K12").Value = 214.63637
$ws.Range("M12").Value = -44.63637

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2569.6
$ws.Range("I18").Value = 2569.6
$ws.Range("K18").Value = 2569.6
$ws.Range("M18").Value = -2285.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 460.375
$ws.Range("I33").Value = 475.4
$ws.Range("J33").Value = 235
$ws.Range("K33").Value = 475.4
$ws.Range("L33").Value = 235
$ws.Range("M33").Value = -246.4
$ws.Range("N33").Value = -693

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1222700
$ws.Range("I98").Value = 4500
$ws.Range("K98").Value = 4500
$ws.Range("M98").Value = -3002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 10308.5625
$ws.Range("I116").Value = 12760.2
$ws.Range("J116").Value = 9194.182000000001
$ws.Range("K116").Value = 12760.2
$ws.Range("L116").Value = 9194.182000000001
$ws.Range("M116").Value = -9318.200000000001
$ws.Range("N116").Value = -16078.182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1222700
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1356.3518
$ws.Range("I135").Value = 1145.8085
$ws.Range("J135").Value = 2770
$ws.Range("K135").Value = 10312.2765
$ws.Range("L135").Value = 24930
$ws.Range("M135").Value = -7777.2765
$ws.Range("N135").Value = -30000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8571.462
$ws.Range("I61").Value = 4261.5
$ws.Range("K61").Value = 4261.5
$ws.Range("M61").Value = -4049.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7754346
$ws.Range("I74").Value = 8548836
$ws.Range("J74").Value = 8072.25
$ws.Range("K74").Value = 8548836
$ws.Range("L74").Value = 8072.25
$ws.Range("M74").Value = -8547962
$ws.Range("N74").Value = -9820.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7754346
$ws.Range("I77").Value = 8548836
$ws.Range("J77").Value = 8072.25
$ws.Range("K77").Value = 42744180
$ws.Range("L77").Value = 40361.25
$ws.Range("M77").Value = -42739812
$ws.Range("N77").Value = -49097.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8571.462
$ws.Range("I136").Value = 4261.5
$ws.Range("K136").Value = 12784.5
$ws.Range("M136").Value = -10234.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1007.6667
$ws.Range("I22").Value = 833.3333
$ws.Range("J22").Value = 1182
$ws.Range("K22").Value = 833.3333
$ws.Range("L22").Value = 1182
$ws.Range("M22").Value = -660.3333
$ws.Range("N22").Value = -1528

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1066.1177
$ws.Range("I134").Value = 1049.5151
$ws.Range("K134").Value = 3148.5453
$ws.Range("M134").Value = -613.5453000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1677.091
$ws.Range("I22").Value = 305
$ws.Range("K22").Value = 305
$ws.Range("M22").Value = 45

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30363.365
$ws.Range("I31").Value = 3531.913
$ws.Range("J31").Value = 64648
$ws.Range("K31").Value = 3531.913
$ws.Range("L31").Value = 64648
$ws.Range("M31").Value = -3236.913
$ws.Range("N31").Value = -65238

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 30363.365
$ws.Range("I34").Value = 3531.913
$ws.Range("J34").Value = 64648
$ws.Range("K34").Value = 3531.913
$ws.Range("L34").Value = 64648
$ws.Range("M34").Value = -3329.913
$ws.Range("N34").Value = -65052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2473.9167
$ws.Range("I99").Value = 1697.6666
$ws.Range("J99").Value = 2732.6667
$ws.Range("K99").Value = 1697.6666
$ws.Range("L99").Value = 2732.6667
$ws.Range("M99").Value = -199.6666
$ws.Range("N99").Value = -5728.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3926.0286
$ws.Range("I122").Value = 2587.611
$ws.Range("J122").Value = 5343.1763
$ws.Range("K122").Value = 7762.833
$ws.Range("L122").Value = 16029.5289
$ws.Range("M122").Value = -5312.833
$ws.Range("N122").Value = -20929.5289

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2473.9167
$ws.Range("I126").Value = 1697.6666
$ws.Range("J126").Value = 2732.6667
$ws.Range("K126").Value = 5092.9998
$ws.Range("L126").Value = 8198.000100000001
$ws.Range("M126").Value = -2622.9998
$ws.Range("N126").Value = -13138.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2716.5881
$ws.Range("I134").Value = 1580.4348
$ws.Range("J134").Value = 5092.1816
$ws.Range("K134").Value = 4741.3044
$ws.Range("L134").Value = 15276.5448
$ws.Range("M134").Value = -2206.3044
$ws.Range("N134").Value = -20346.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6794.816
$ws.Range("I7").Value = 5193.7666
$ws.Range("J7").Value = 12798.75
$ws.Range("K7").Value = 5193.7666
$ws.Range("L7").Value = 12798.75
$ws.Range("M7").Value = -5081.7666
$ws.Range("N7").Value = -13022.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6253.933
$ws.Range("I40").Value = 5831.548
$ws.Range("J40").Value = 12167.333
$ws.Range("K40").Value = 5831.548
$ws.Range("L40").Value = 5831.548
$ws.Range("M40").Value = -5695.548
$ws.Range("N40").Value = -12439.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6794.816
$ws.Range("I126").Value = 5193.7666
$ws.Range("J126").Value = 12798.75
$ws.Range("K126").Value = 15581.2998
$ws.Range("L126").Value = 38396.25
$ws.Range("M126").Value = -13111.2998
$ws.Range("N126").Value = -43336.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4095.4546
$ws.Range("I132").Value = 3397.4102
$ws.Range("J132").Value = 9540.200000000001
$ws.Range("K132").Value = 10192.2306
$ws.Range("L132").Value = 28620.6
$ws.Range("M132").Value = -7662.230599999999
$ws.Range("N132").Value = -33680.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3948.5854
$ws.Range("I136").Value = 1922.1111
$ws.Range("J136").Value = 18539.2
$ws.Range("K136").Value = 5766.3333
$ws.Range("L136").Value = 55617.60000000001
$ws.Range("M136").Value = -3216.3333
$ws.Range("N136").Value = -60717.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2115.5938
$ws.Range("I126").Value = 1707.5652
$ws.Range("K126").Value = 5122.6956
$ws.Range("M126").Value = -2652.6956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2346.1404
$ws.Range("I136").Value = 1737.9584
$ws.Range("J136").Value = 5589.778
$ws.Range("K136").Value = 5213.8752
$ws.Range("L136").Value = 16769.334
$ws.Range("M136").Value = -2663.8752
